$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: update the date (was 44874, now 44876) and add the missing Video (column C) ---
$ws.Range("A28").Value = 44876
$ws.Range("C28").Value = "Ground Reference Manuevers"

# --- Helper: copy the date formatting from A28 down to a new row's A cell, then set its value ---
function Set-DateRow($row, $dateSerial, $video) {
    $ws.Range("A28").Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$row").Value = $dateSerial
    $ws.Range("B$row").Value = "Practicing Landings"
    $ws.Range("C$row").Value = $video
}

Set-DateRow 29 44876 "Closer Look: Taxi Tips"
Set-DateRow 30 44876 "Engines"
Set-DateRow 31 44876 "Air Facts: Engine Suspicion"
Set-DateRow 32 44878 "Aerodynamics"
Set-DateRow 33 44878 "Air Closer Look: Angle of Attack"
Set-DateRow 34 44878 "Slow Flight"
Set-DateRow 35 44878 "Closer Look: Change of Scenery"
Set-DateRow 36 44976 "Stalls"
Set-DateRow 37 44978 "Air Facts:Stall Rhetoric"

# --- Rows 38-39: the date was mistyped ("2/21/20223") so Excel stored it as plain text ---
$ws.Range("A38").Value = "2/21/20223"
$ws.Range("B38").Value = "Practicing Landings"
$ws.Range("C38").Value = "Normal Landings"

$ws.Range("A39").Value = "2/21/20223"
$ws.Range("B39").Value = "Practicing Landings"
$ws.Range("C39").Value = "Air Facts: Down to Earth"

# --- Column A widened (auto-fit) to accommodate the new, longer dates ---
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(1).ColumnWidth = 9.6

# --- Restore view: scroll down near the bottom of the newly entered data ---
$ws.Range("K54").Select() | Out-Null
